$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where both Price (D) and Volume(1h) (E) columns change.
# The "'" prefix on D forces Excel to store the value as text even when it
# looks numeric (e.g. "307.14"); resetting Style back to "Normal" afterwards
# removes the quote-prefix formatting flag so the cell keeps the default style.
$priceVolumeUpdates = @(
    @{ Row = 2;  D = "24.337.31";    E = "  +8.85%  " },
    @{ Row = 3;  D = "1.677.56";     E = "  +4.67%  " },
    @{ Row = 4;  D = "1.004";        E = "  -0.37%  " },
    @{ Row = 5;  D = "307.14";       E = "  +6.29%  " },
    @{ Row = 6;  D = "0.9977";       E = "  +0.19%  " },
    @{ Row = 7;  D = "0.3710";       E = "  -0.16%  " },
    @{ Row = 8;  D = "0.3442";       E = "  +1.65%  " },
    @{ Row = 9;  D = "48.23";        E = "  +13.22%  " },
    @{ Row = 10; D = "1.183";        E = "  +3.29%  " },
    @{ Row = 11; D = "0.07249";      E = "  +2.51%  " },
    @{ Row = 12; D = "1.001";        E = "  -0.15%  " },
    @{ Row = 13; D = "20.36";        E = "  +2.36%  " },
    @{ Row = 14; D = "6.103";        E = "  +2.75%  " },
    @{ Row = 16; D = "1.677.46";     E = "  +4.69%  " },
    @{ Row = 17; D = "0.00001110";   E = "  +1.92%  " },
    @{ Row = 18; D = "0.9983";       E = "  +0.27%  " },
    @{ Row = 19; D = "0.06720";      E = "  +1.12%  " },
    @{ Row = 20; D = "81.20";        E = "  +3.45%  " },
    @{ Row = 21; D = "16.45";        E = "  +1.22%  " },
    @{ Row = 22; D = "6.093";        E = "  +0.61%  " },
    @{ Row = 23; D = "11.95";        E = "  +0.94%  " },
    @{ Row = 24; D = "24.318.55";    E = "  +8.59%  " },
    @{ Row = 25; D = "2.430";        E = "  +1.89%  " },
    @{ Row = 27; D = "2.658";        E = "  +6.11%  " },
    @{ Row = 28; D = "152.11";       E = "  +0.66%  " },
    @{ Row = 31; D = "127.28";       E = "  +5.15%  " },
    @{ Row = 32; D = "6.318";        E = "  +5.00%  " },
    @{ Row = 33; D = "4.030";        E = "  -4.50%  " },
    @{ Row = 34; D = "0.9667";       E = "  +1.58%  " },
    @{ Row = 35; D = "1.746";        E = "  +8.30%  " },
    @{ Row = 36; D = "0.08468";      E = "  +2.33%  " },
    @{ Row = 37; D = "8.982";        E = "  +3.31%  " },
    @{ Row = 38; D = "12.28";        E = "  +3.88%  " },
    @{ Row = 39; D = "0.06415";      E = "  +3.83%  " },
    @{ Row = 40; D = "5.332";        E = "  -0.08%  " },
    @{ Row = 42; D = "1.264";        E = "  +1.76%  " },
    @{ Row = 44; D = "0.6172";       E = "  +3.46%  " },
    @{ Row = 45; D = "0.9975";       E = "  +0.32%  " },
    @{ Row = 49; D = "126.95";       E = "  +0.91%  " },
    @{ Row = 50; D = "2.024";        E = "  +2.33%  " },
    @{ Row = 51; D = "0.07209";      E = "  +5.21%  " }
)

foreach ($u in $priceVolumeUpdates) {
    $dCell = $ws.Cells.Item($u.Row, 4)
    $dCell.Value = "'" + $u.D
    $dCell.Style = "Normal"
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

# Rows where only Volume(1h) (E) changes
$volumeOnlyUpdates = @(
    @{ Row = 15; E = "  +1.01%  " },
    @{ Row = 26; E = "  -11.94%  " },
    @{ Row = 29; E = "  -1.00%  " },
    @{ Row = 30; E = "  +4.52%  " },
    @{ Row = 41; E = "  +5.10%  " },
    @{ Row = 43; E = "  +4.07%  " }
)

foreach ($u in $volumeOnlyUpdates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

# Row 46: only Price (D) changes
$d46 = $ws.Cells.Item(46, 4)
$d46.Value = "'3.775"
$d46.Style = "Normal"

# Rows 47 and 48: Decentraland and EnergySwap swap positions (with updated data)
$ws.Cells.Item(47, 2).Value = "Decentraland"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$d47 = $ws.Cells.Item(47, 4)
$d47.Value = "'0.5939"
$d47.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +3.41%  "

$ws.Cells.Item(48, 2).Value = "EnergySwap"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$d48 = $ws.Cells.Item(48, 4)
$d48.Value = "'13.02"
$d48.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -1.46%  "
